# Apply the diff to the DeliverySequence workbook.
# The workbook has 6 sheets named "z1,1", "z2,1", "z1,2", "z2,2", "z1,3", "z2,3"
# corresponding to sheet1.xml .. sheet6.xml. Only sheets z2,1 / z1,2 / z1,3 / z2,3
# change in this edit.

$wb = $excel.ActiveWorkbook

# --- Sheet "z2,1" (sheet2.xml) ---
$ws = $wb.Worksheets.Item("z2,1")
$ws.Range("F1").Value = 1
$ws.Range("G1").Value = 0
$ws.Range("D2").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0

# --- Sheet "z1,2" (sheet3.xml) ---
$ws = $wb.Worksheets.Item("z1,2")
$ws.Range("B1").Value = 0
$ws.Range("B6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("B10").Value = 0

# --- Sheet "z1,3" (sheet5.xml) ---
$ws = $wb.Worksheets.Item("z1,3")
$ws.Range("G8").Value = 1

# --- Sheet "z2,3" (sheet6.xml) ---
$ws = $wb.Worksheets.Item("z2,3")
$ws.Range("F1").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("I4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("E6").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("F10").Value = 0
